# Wheat-price-data.xlsx — "reorganized things. Fixed a thing"
#
# 1. C2 held the literal text "-" (a placeholder for the first month, which
#    has no prior month to compute a percent change against). It's replaced
#    with an actual numeric value (0.5) so the column is fully numeric.
#    Excel drops the now-unreferenced shared string "-" from the shared
#    string table on save, which cascades: every later shared-string index
#    (all the other B-column price labels, and the C1 "Percent Change"
#    header) shifts down by one automatically — that part requires no
#    explicit action here.
# 2. The percent-change column (C2:C211) gets an explicit "0.000" number
#    format (a fresh custom numFmt) instead of the default General format.
# 3. The sheet's saved view is reset: no more frozen/scrolled
#    "topLeftCell", and the active selection moves to E5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Replace the "-" placeholder in C2 with a real number.
$ws.Range("C2").Value = 0.5

# 2. Apply the new numeric format to the whole percent-change column.
$ws.Range("C2:C211").NumberFormat = "0.000"

# 3. Update the sheet view / selection.
$ws.Range("E5").Select()
